# Update attendance counts in Sheet1 (final code with comment - tut 6)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1

# Row 5
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Row 6
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

# Rows 7-18: mark Absent column (H) as 1
for ($r = 7; $r -le 18; $r++) {
    $ws.Cells.Item($r, 8).Value = 1
}
